$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 247.25
$ws.Range("I33").Value = 196.8
$ws.Range("K33").Value = 196.8
$ws.Range("M33").Value = 32.19999999999999
$ws.Range("H43").Value = 15598.25
$ws.Range("I43").Value = 22979.6
$ws.Range("J43").Value = 3296
$ws.Range("K43").Value = 22979.6
$ws.Range("L43").Value = 3296
$ws.Range("M43").Value = -22910.6
$ws.Range("N43").Value = -3434
$ws.Range("H75").Value = 100000
$ws.Range("J75").Value = 100000
$ws.Range("L75").Value = 100000
$ws.Range("N75").Value = -101872
$ws.Range("H78").Value = 100000
$ws.Range("J78").Value = 100000
$ws.Range("L78").Value = 300000
$ws.Range("N78").Value = -309360
$ws.Range("H98").Value = 1547.9131
$ws.Range("I98").Value = 1671.2778
$ws.Range("K98").Value = 1671.2778
$ws.Range("M98").Value = -173.2778000000001
$ws.Range("H106").Value = 15838.059
$ws.Range("I106").Value = 7649.4
$ws.Range("K106").Value = 7649.4
$ws.Range("M106").Value = -7018.4
$ws.Range("H107").Value = 6192.3
$ws.Range("I107").Value = 7690.375
$ws.Range("K107").Value = 7690.375
$ws.Range("M107").Value = -5770.375
$ws.Range("H116").Value = 39942.855
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 45266.668
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 45266.668
$ws.Range("M116").Value = -4558
$ws.Range("N116").Value = -52150.668
$ws.Range("H122").Value = 1547.9131
$ws.Range("I122").Value = 1671.2778
$ws.Range("K122").Value = 5013.8334
$ws.Range("M122").Value = -2563.8334
$ws.Range("H138").Value = 2172.7163
$ws.Range("I138").Value = 1201
$ws.Range("J138").Value = 2440.776
$ws.Range("K138").Value = 3603
$ws.Range("L138").Value = 7322.328
$ws.Range("M138").Value = 1537
$ws.Range("N138").Value = -17602.328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 201297.4
$ws.Range("I32").Value = 209027.69
$ws.Range("K32").Value = 209027.69
$ws.Range("M32").Value = -208740.69
$ws.Range("H61").Value = 4082.25
$ws.Range("I61").Value = 3759.7693
$ws.Range("K61").Value = 3759.7693
$ws.Range("M61").Value = -3547.7693
$ws.Range("H110").Value = 1525
$ws.Range("I110").Value = 1700
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1700
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 345
$ws.Range("N110").Value = -5090
$ws.Range("H136").Value = 4082.25
$ws.Range("I136").Value = 3759.7693
$ws.Range("K136").Value = 11279.3079
$ws.Range("M136").Value = -8729.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 363.66666
$ws.Range("I22").Value = 195.5
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 195.5
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -22.5
$ws.Range("N22").Value = -1046
$ws.Range("H105").Value = 2720.4412
$ws.Range("I105").Value = 1849.875
$ws.Range("K105").Value = 1849.875
$ws.Range("M105").Value = -102.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 10679.286
$ws.Range("J57").Value = 15000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16120
$ws.Range("H93").Value = 7883.8335
$ws.Range("I93").Value = 7883.8335
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 7883.8335
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -6011.8335
$ws.Range("N93").ClearContents()
$ws.Range("H99").Value = 3832.8333
$ws.Range("I99").Value = 3499
$ws.Range("J99").Value = 3999.75
$ws.Range("K99").Value = 3499
$ws.Range("L99").Value = 3999.75
$ws.Range("M99").Value = -2001
$ws.Range("N99").Value = -6995.75
$ws.Range("H126").Value = 3832.8333
$ws.Range("I126").Value = 3499
$ws.Range("J126").Value = 3999.75
$ws.Range("K126").Value = 10497
$ws.Range("L126").Value = 11999.25
$ws.Range("M126").Value = -8027
$ws.Range("N126").Value = -16939.25
$ws.Range("H139").Value = 43926.668
$ws.Range("J139").Value = 43926.668
$ws.Range("L139").Value = 43926.668
$ws.Range("N139").Value = -54206.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2619.5
$ws.Range("I12").Value = 2517.5
$ws.Range("K12").Value = 7552.5
$ws.Range("M12").Value = -7379.5
$ws.Range("H38").Value = 66.7
$ws.Range("I38").Value = 48.266666
$ws.Range("J38").Value = 122
$ws.Range("K38").Value = 144.799998
$ws.Range("L38").Value = 366
$ws.Range("M38").Value = 202.200002
$ws.Range("N38").Value = -1060
$ws.Range("H86").Value = 399.5
$ws.Range("J86").Value = 399.5
$ws.Range("L86").Value = 1198.5
$ws.Range("N86").Value = -3570.5
$ws.Range("H87").Value = 9400
$ws.Range("I87").Value = 6880
$ws.Range("K87").Value = 20640
$ws.Range("M87").Value = -19392
$ws.Range("H89").Value = 399.5
$ws.Range("J89").Value = 399.5
$ws.Range("L89").Value = 3595.5
$ws.Range("N89").Value = -15451.5
$ws.Range("H90").Value = 9400
$ws.Range("I90").Value = 6880
$ws.Range("K90").Value = 61920
$ws.Range("M90").Value = -55680
$ws.Range("H92").Value = 596.6667
$ws.Range("I92").Value = 596.6667
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1790.0001
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -542.0001
$ws.Range("N92").ClearContents()
$ws.Range("H103").Value = 293.66666
$ws.Range("J103").Value = 539
$ws.Range("L103").Value = 1617
$ws.Range("N103").Value = -3375
$ws.Range("H131").Value = 2446.9534
$ws.Range("J131").Value = 2446.9534
$ws.Range("L131").Value = 7340.860199999999
$ws.Range("N131").Value = -17420.8602

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2096.25
$ws.Range("I126").Value = 1929.2858
$ws.Range("J126").Value = 2485.8333
$ws.Range("K126").Value = 5787.857400000001
$ws.Range("L126").Value = 7457.499899999999
$ws.Range("M126").Value = -3317.857400000001
$ws.Range("N126").Value = -12397.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 3481.647
$ws.Range("I7").Value = 2808.9092
$ws.Range("K7").Value = 2808.9092
$ws.Range("M7").Value = -2696.9092
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H35").Value = 13372.125
$ws.Range("I35").Value = 2823.6667
$ws.Range("J35").Value = 45017.5
$ws.Range("K35").Value = 2823.6667
$ws.Range("L35").Value = 45017.5
$ws.Range("M35").Value = -2487.6667
$ws.Range("N35").Value = -45689.5
$ws.Range("H40").Value = 3580.9092
$ws.Range("I40").Value = 2180
$ws.Range("J40").Value = 4748.3335
$ws.Range("K40").Value = 2180
$ws.Range("L40").Value = 4748.3335
$ws.Range("M40").Value = -2044
$ws.Range("N40").Value = -5020.3335
$ws.Range("H46").Value = 4276.6523
$ws.Range("I46").Value = 1196.5
$ws.Range("J46").Value = 4570
$ws.Range("K46").Value = 1196.5
$ws.Range("L46").Value = 4570
$ws.Range("M46").Value = -1008.5
$ws.Range("N46").Value = -4946
$ws.Range("H126").Value = 3481.647
$ws.Range("I126").Value = 2808.9092
$ws.Range("K126").Value = 8426.7276
$ws.Range("M126").Value = -5956.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3180
$ws.Range("I14").Value = 3000
$ws.Range("J14").Value = 3450
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3450
$ws.Range("M14").Value = -2832
$ws.Range("N14").Value = -3786
$ws.Range("H32").Value = 15000
$ws.Range("I32").Value = 15000
$ws.Range("K32").Value = 15000
$ws.Range("M32").Value = -14683
$ws.Range("H37").Value = 5029
$ws.Range("J37").Value = 5029
$ws.Range("L37").Value = 5029
$ws.Range("N37").Value = -5435
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H107").Value = 838.6905
$ws.Range("I107").Value = 785.69446
$ws.Range("J107").Value = 1156.6666
$ws.Range("K107").Value = 2357.08338
$ws.Range("L107").Value = 3469.9998
$ws.Range("M107").Value = -437.08338
$ws.Range("N107").Value = -7309.9998
$ws.Range("H122").Value = 77299.8
$ws.Range("I122").Value = 2472.875
$ws.Range("K122").Value = 7418.625
$ws.Range("M122").Value = -4968.625
$ws.Range("H126").Value = 2132
$ws.Range("I126").Value = 2116.2666
$ws.Range("K126").Value = 6348.7998
$ws.Range("M126").Value = -3878.7998
$ws.Range("H136").Value = 1969.75
$ws.Range("I136").Value = 2063.8
$ws.Range("K136").Value = 6191.400000000001
$ws.Range("M136").Value = -3641.400000000001

Write-Output "Applied changes to all sheets"